# Update Name of Algo
# Apply updated numeric values produced by the new algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E7").Value  = 13.377
$ws.Range("A10").Value = -20.926
$ws.Range("A12").Value = -21.694
$ws.Range("B13").Value = 6.606
$ws.Range("A18").Value = -21.694
$ws.Range("E20").Value = 12.932

$wb.Save()
